$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(109, 8).Value = 48500
$ws.Cells.Item(109, 10).Value = 48500
$ws.Cells.Item(109, 12).Value = 48500
$ws.Cells.Item(109, 14).Value = -51274
$ws.Cells.Item(116, 8).Value = 11371785
$ws.Cells.Item(116, 9).Value = 19235406
$ws.Cells.Item(116, 11).Value = 19235406
$ws.Cells.Item(116, 13).Value = -19231964
$ws.Cells.Item(132, 8).Value = 17371.541
$ws.Cells.Item(132, 9).Value = 18531.877
$ws.Cells.Item(132, 11).Value = 55595.631
$ws.Cells.Item(132, 13).Value = -53065.631
$ws.Cells.Item(138, 8).Value = 1013649.44
$ws.Cells.Item(138, 9).Value = 2205.359
$ws.Cells.Item(138, 10).Value = 1671088.1
$ws.Cells.Item(138, 11).Value = 6616.076999999999
$ws.Cells.Item(138, 12).Value = 5013264.300000001
$ws.Cells.Item(138, 13).Value = -1476.076999999999
$ws.Cells.Item(138, 14).Value = -5023544.300000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 15153295
$ws.Cells.Item(2, 10).Value = 66671590
$ws.Cells.Item(2, 12).Value = 66671590
$ws.Cells.Item(2, 14).Value = -66671816
$ws.Cells.Item(32, 8).Value = 2574683.2
$ws.Cells.Item(32, 9).Value = 3082109.8
$ws.Cells.Item(32, 11).Value = 3082109.8
$ws.Cells.Item(32, 13).Value = -3081822.8
$ws.Cells.Item(45, 8).Value = 4005.1875
$ws.Cells.Item(45, 9).Value = 2950.8
$ws.Cells.Item(45, 11).Value = 2950.8
$ws.Cells.Item(45, 13).Value = -2573.8
$ws.Cells.Item(61, 8).Value = 24393162
$ws.Cells.Item(61, 9).Value = 2102.3845
$ws.Cells.Item(61, 11).Value = 2102.3845
$ws.Cells.Item(61, 13).Value = -1890.3845
$ws.Cells.Item(74, 8).Value = 43662.68
$ws.Cells.Item(74, 9).Value = 79013.46000000001
$ws.Cells.Item(74, 10).Value = 5366
$ws.Cells.Item(74, 11).Value = 79013.46000000001
$ws.Cells.Item(74, 12).Value = 5366
$ws.Cells.Item(74, 13).Value = -78139.46000000001
$ws.Cells.Item(74, 14).Value = -7114
$ws.Cells.Item(77, 8).Value = 43662.68
$ws.Cells.Item(77, 9).Value = 79013.46000000001
$ws.Cells.Item(77, 10).Value = 5366
$ws.Cells.Item(77, 11).Value = 395067.3
$ws.Cells.Item(77, 12).Value = 26830
$ws.Cells.Item(77, 13).Value = -390699.3
$ws.Cells.Item(77, 14).Value = -35566
$ws.Cells.Item(116, 8).Value = 15153295
$ws.Cells.Item(116, 10).Value = 66671590
$ws.Cells.Item(116, 12).Value = 66671590
$ws.Cells.Item(116, 14).Value = -66676178
$ws.Cells.Item(136, 8).Value = 24393162
$ws.Cells.Item(136, 9).Value = 2102.3845
$ws.Cells.Item(136, 11).Value = 6307.1535
$ws.Cells.Item(136, 13).Value = -3757.1535

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 15153295
$ws.Cells.Item(3, 10).Value = 66671590
$ws.Cells.Item(3, 12).Value = 66671590
$ws.Cells.Item(3, 14).Value = -66671818
$ws.Cells.Item(86, 8).Value = 9299019
$ws.Cells.Item(86, 10).Value = 5488
$ws.Cells.Item(86, 12).Value = 5488
$ws.Cells.Item(86, 14).Value = -7734
$ws.Cells.Item(89, 8).Value = 9299019
$ws.Cells.Item(89, 10).Value = 5488
$ws.Cells.Item(89, 12).Value = 27440
$ws.Cells.Item(89, 14).Value = -38672
$ws.Cells.Item(105, 8).Value = 3326.923
$ws.Cells.Item(105, 9).Value = 2551.45
$ws.Cells.Item(105, 10).Value = 4143.2104
$ws.Cells.Item(105, 11).Value = 2551.45
$ws.Cells.Item(105, 12).Value = 4143.2104
$ws.Cells.Item(105, 13).Value = -804.4499999999998
$ws.Cells.Item(105, 14).Value = -7637.2104
$ws.Cells.Item(134, 8).Value = 5322332
$ws.Cells.Item(134, 9).Value = 8334155.5
$ws.Cells.Item(134, 11).Value = 25002466.5
$ws.Cells.Item(134, 13).Value = -24999931.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5506.6665
$ws.Cells.Item(31, 9).Value = 1995.7142
$ws.Cells.Item(31, 10).Value = 7092.2583
$ws.Cells.Item(31, 11).Value = 1995.7142
$ws.Cells.Item(31, 12).Value = 7092.2583
$ws.Cells.Item(31, 13).Value = -1700.7142
$ws.Cells.Item(31, 14).Value = -7682.2583
$ws.Cells.Item(34, 8).Value = 5506.6665
$ws.Cells.Item(34, 9).Value = 1995.7142
$ws.Cells.Item(34, 10).Value = 7092.2583
$ws.Cells.Item(34, 11).Value = 1995.7142
$ws.Cells.Item(34, 12).Value = 7092.2583
$ws.Cells.Item(34, 13).Value = -1793.7142
$ws.Cells.Item(34, 14).Value = -7496.2583
$ws.Cells.Item(62, 8).Value = 5646.593
$ws.Cells.Item(62, 9).Value = 5266.6
$ws.Cells.Item(62, 10).Value = 5870.1177
$ws.Cells.Item(62, 11).Value = 5266.6
$ws.Cells.Item(62, 12).Value = 5870.1177
$ws.Cells.Item(62, 13).Value = -4642.6
$ws.Cells.Item(62, 14).Value = -7118.1177
$ws.Cells.Item(65, 8).Value = 5646.593
$ws.Cells.Item(65, 9).Value = 5266.6
$ws.Cells.Item(65, 10).Value = 5870.1177
$ws.Cells.Item(65, 11).Value = 26333
$ws.Cells.Item(65, 12).Value = 29350.5885
$ws.Cells.Item(65, 13).Value = -23213
$ws.Cells.Item(65, 14).Value = -35590.5885
$ws.Cells.Item(132, 8).Value = 3867.125
$ws.Cells.Item(132, 9).Value = 2635.5356
$ws.Cells.Item(132, 10).Value = 5591.35
$ws.Cells.Item(132, 11).Value = 7906.6068
$ws.Cells.Item(132, 12).Value = 16774.05
$ws.Cells.Item(132, 13).Value = -5376.6068
$ws.Cells.Item(132, 14).Value = -21834.05
$ws.Cells.Item(134, 8).Value = 3525.1892
$ws.Cells.Item(134, 9).Value = 2393.3809
$ws.Cells.Item(134, 10).Value = 5010.6875
$ws.Cells.Item(134, 11).Value = 7180.1427
$ws.Cells.Item(134, 12).Value = 15032.0625
$ws.Cells.Item(134, 13).Value = -4645.1427
$ws.Cells.Item(134, 14).Value = -20102.0625
$ws.Cells.Item(141, 8).Value = 62165.54
$ws.Cells.Item(141, 10).Value = 62165.54
$ws.Cells.Item(141, 12).Value = 62165.54
$ws.Cells.Item(141, 14).Value = -72525.54000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1223.6774
$ws.Cells.Item(5, 9).Value = 774.3889
$ws.Cells.Item(5, 11).Value = 2323.1667
$ws.Cells.Item(5, 13).Value = -2211.1667
$ws.Cells.Item(8, 8).Value = 333.8
$ws.Cells.Item(8, 9).Value = 333.8
$ws.Cells.Item(8, 11).Value = 1001.4
$ws.Cells.Item(8, 13).Value = -862.4000000000001
$ws.Cells.Item(131, 8).Value = 49355.668
$ws.Cells.Item(131, 9).Value = 865
$ws.Cells.Item(131, 10).Value = 54459.95
$ws.Cells.Item(131, 11).Value = 2595
$ws.Cells.Item(131, 12).Value = 163379.85
$ws.Cells.Item(131, 13).Value = 2445
$ws.Cells.Item(131, 14).Value = -173459.85
$ws.Cells.Item(135, 8).Value = 1223.6774
$ws.Cells.Item(135, 9).Value = 774.3889
$ws.Cells.Item(135, 11).Value = 6969.5001
$ws.Cells.Item(135, 13).Value = -4434.5001
$ws.Cells.Item(136, 8).Value = 17859238
$ws.Cells.Item(136, 9).Value = 17859238
$ws.Cells.Item(136, 11).Value = 53577714
$ws.Cells.Item(136, 13).Value = -53572614
$ws.Cells.Item(139, 8).Value = 39703.31
$ws.Cells.Item(139, 9).Value = 55442.527
$ws.Cells.Item(139, 11).Value = 166327.581
$ws.Cells.Item(139, 13).Value = -161187.581

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1545.2858
$ws.Cells.Item(22, 10).Value = 2358.75
$ws.Cells.Item(22, 12).Value = 2358.75
$ws.Cells.Item(22, 14).Value = -2948.75
$ws.Cells.Item(27, 8).Value = 1545.2858
$ws.Cells.Item(27, 10).Value = 2358.75
$ws.Cells.Item(27, 12).Value = 2358.75
$ws.Cells.Item(27, 14).Value = -2572.75
$ws.Cells.Item(61, 8).Value = 3229770
$ws.Cells.Item(61, 9).Value = 4764504
$ws.Cells.Item(61, 11).Value = 4764504
$ws.Cells.Item(61, 13).Value = -4764302
$ws.Cells.Item(113, 8).Value = 3229770
$ws.Cells.Item(113, 9).Value = 4764504
$ws.Cells.Item(113, 11).Value = 4764504
$ws.Cells.Item(113, 13).Value = -4762334
$ws.Cells.Item(132, 8).Value = 8200610
$ws.Cells.Item(132, 9).Value = 12197679
$ws.Cells.Item(132, 10).Value = 6617.85
$ws.Cells.Item(132, 11).Value = 36593037
$ws.Cells.Item(132, 12).Value = 19853.55
$ws.Cells.Item(132, 13).Value = -36590507
$ws.Cells.Item(132, 14).Value = -24913.55
$ws.Cells.Item(136, 8).Value = 7627.5737
$ws.Cells.Item(136, 9).Value = 2483.4443
$ws.Cells.Item(136, 11).Value = 7450.3329
$ws.Cells.Item(136, 13).Value = -4900.3329

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 180000000
$ws.Cells.Item(5, 9).Value = 180000000
$ws.Cells.Item(5, 11).Value = 180000000
$ws.Cells.Item(5, 13).Value = -179999888
$ws.Cells.Item(62, 8).Value = 48639.855
$ws.Cells.Item(62, 9).Value = 56403.168
$ws.Cells.Item(62, 10).Value = 2060
$ws.Cells.Item(62, 11).Value = 56403.168
$ws.Cells.Item(62, 12).Value = 2060
$ws.Cells.Item(62, 13).Value = -55779.168
$ws.Cells.Item(62, 14).Value = -3308
$ws.Cells.Item(65, 8).Value = 48639.855
$ws.Cells.Item(65, 9).Value = 56403.168
$ws.Cells.Item(65, 10).Value = 2060
$ws.Cells.Item(65, 11).Value = 282015.84
$ws.Cells.Item(65, 12).Value = 10300
$ws.Cells.Item(65, 13).Value = -278895.84
$ws.Cells.Item(65, 14).Value = -16540
$ws.Cells.Item(113, 8).Value = 8148.795
$ws.Cells.Item(113, 9).Value = 13394.333
$ws.Cells.Item(113, 11).Value = 40182.999
$ws.Cells.Item(113, 13).Value = -38012.999
$ws.Cells.Item(126, 8).Value = 1006.5
$ws.Cells.Item(126, 9).Value = 927.5
$ws.Cells.Item(126, 11).Value = 2782.5
$ws.Cells.Item(126, 13).Value = -312.5
$ws.Cells.Item(132, 8).Value = 4580.2666
$ws.Cells.Item(132, 9).Value = 4441.2812
$ws.Cells.Item(132, 10).Value = 4922.385
$ws.Cells.Item(132, 11).Value = 13323.8436
$ws.Cells.Item(132, 12).Value = 14767.155
$ws.Cells.Item(132, 13).Value = -10793.8436
$ws.Cells.Item(132, 14).Value = -19827.155
$ws.Cells.Item(136, 8).Value = 17722526
$ws.Cells.Item(136, 9).Value = 24391568
$ws.Cells.Item(136, 11).Value = 73174704
$ws.Cells.Item(136, 13).Value = -73172154
